$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header labels: "_old" -> "_FV2210" and "_new" -> "_FV2304"
#    (these suffixes only occur in the ten header cells of each block, on row 1)
$ws.Cells.Replace("_old", "_FV2210")
$ws.Cells.Replace("_new", "_FV2304")

# 2. Turn the used range into an Excel Table ("Table1") so the renamed headers
#    become the table's column headers.
$dataRange = $ws.Range("A1:U53")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# 3. Freeze the header row (row 1) on the worksheet.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
